# Generate Report for Handback
#
# 1. Update the "Ready for handoff" status text (shared across Overview,
#    zh-cn and de-de sheets) to "Handback transform failed".
# 2. Fill in the "Error Detail" column (column P) for the 9464d5ae... row
#    on the zh-cn and de-de sheets with the handback/handoff file-name
#    mismatch message (language specific).
# 3. Widen column P ("Error Detail") on the zh-cn and de-de sheets to a
#    stored width of 40.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# --- Status text update (was "Ready for handoff") ---
$ws_overview.Range("E3").Value = $newStatus
$ws_overview.Range("F3").Value = $newStatus
$ws_zhcn.Range("C3").Value = $newStatus
$ws_dede.Range("C3").Value = $newStatus

# --- Error Detail messages for the 9464d5ae... row (row 3) ---
$zhcnError = "Handback file name: 5gi2w5xt.ihn is different with handoff file name: 9464d5ae-7490-4400-a33f-0c01fcd101f2.5ed681592b9fa81b779d348c336ef8bef2494708.zh-cn."
$dedeError = "Handback file name: 5gi2w5xt.ihn is different with handoff file name: 9464d5ae-7490-4400-a33f-0c01fcd101f2.5ed681592b9fa81b779d348c336ef8bef2494708.de-de."

$ws_zhcn.Range("P3").Value = $zhcnError
$ws_dede.Range("P3").Value = $dedeError

# --- Widen the Error Detail column (P / column 16) to a stored width of 40 ---
$targetColumnWidth = 39 + (1 / 6)
$ws_zhcn.Columns.Item(16).ColumnWidth = $targetColumnWidth
$ws_dede.Columns.Item(16).ColumnWidth = $targetColumnWidth
